# Product backlog uitgebreid: vul de nog lege User story / Taken / Acceptatiecriteria
# cellen in voor de requirements 7 t/m 11 (rijen 8, 10, 11 en 12) en werk de
# tekst van D9 bij (toegevoegde punt aan het einde van de zin).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product backlog")

# Rij 8 ("7. De miniatuur waterpompen ...") - Taken en Acceptatiecriteria aanvullen
$ws.Range("C8").Value = "1. Overleg met opdrachtgever binnen hoeveel minuten.`n2. Gebaseerd op het overleg berekenen hoeveel water de pompen moeten kunnen verplaatsen om dit te bereiken. `n3. Pomp uitkiezen en bestellen`n4. Pomp implementeren"
$ws.Range("D8").Value = "De pompen kunnen de ballasttanks binnen … minuten leeg en vol pompen/laten lopen. "

# Rij 9 ("8. Het schaalmodel lijkt sprekend ...") - kleine tekstcorrectie (punt toegevoegd)
$ws.Range("D9").Value = "Het schaalmodel is op 1:100 schaal van de Floating Farm, met zo veel mogelijk detail."

# Rij 10 ("9. De Ballasttanks kunnen de Floating Farm ... graden kantelen.")
$ws.Range("B10").Value = "Als eigenaar wil ik dat de Ballasttanks de Floating Farm … graden kunnen kantelen, zodat mijn boerderij altijd waterpas kan staan. "
$ws.Range("C10").Value = "1. Berekenen hoeveel drijfvermogen de Ballasttanks moeten hebben.`n2. Maak de Ballasttanks gebaseerd op deze berekeningen."
$ws.Range("D10").Value = "De ballasttanks hebben het drijfvermogen/inhoud om de Floating Farm … graden te draaien."
$ws.Range("E10").Value = "M"

# Rij 11 ("10. Het complete systeem heeft een maximale afwijking ...")
$ws.Range("B11").Value = "Als eigenaar wil ik dat het systeem niet meer dan … graden afwijkt om te zorgen dat mijn lezingen nauwkeurig en correct zijn. "
$ws.Range("C11").Value = "1. Maximale afwijking berekenen`n2.  Huidige afwijkingen bekijken`n3. Aanpassingen implementeren zo nodig"
$ws.Range("D11").Value = "Het systeem heeft ten alle tijden niet meer dan … graden afwijking tot waterpas."
$ws.Range("E11").Value = "M"

# Rij 12 ("11. Het systeem is bruikbaar voor toekomstige designs ...")
$ws.Range("B12").Value = "Als eigenaar wil ik dat het design uitbreidbaar en herbruikbaar is, om te zorgen dat dit probleem opgelost is voordat het onstaat bij volgende Floating Farms. "
$ws.Range("C12").Value = "1. Onderzoeken wat er nodig is om een design herbruikbaar te maken.`n2. Rekening houden met de onderzochte dingen. "
$ws.Range("D12").Value = "Het systeem is gemakkelijk toe te passen bij`n toekomstige Floating Farms en is ook makkelijk uitbreidbaar. "
$ws.Range("D12").WrapText = $true
$ws.Range("E12").Value = "S"

# Rijhoogtes bijwerken zodat de (nu langere) tekst volledig zichtbaar is
$ws.Rows.Item(8).RowHeight = 72.5
$ws.Rows.Item(10).RowHeight = 43.5
$ws.Rows.Item(11).RowHeight = 43.5
$ws.Rows.Item(12).RowHeight = 43.5

# Scroll-/selectiepositie van het werkblad bijwerken
$null = $ws.Select()
$null = $ws.Range("E12").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
